$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (Resources stays in B, existing
# Time/Space pairs shift right to make room for the new "Notes" and
# "Python Concepts" columns).
$ws.Columns("C:D").Insert()

# New column widths match the rest of the wide text columns (A:B).
$ws.Columns("C:D").ColumnWidth = 39.5

# --- Header row (row 1) ---
$ws.Range("C1").Value = "Notes"
$ws.Range("D1").Value = "Python Concepts"

# --- Row 2 (Fibonacci Number row) ---
$ws.Range("B2").Value = "s"
$ws.Range("B2").Font.Bold = $true

# --- Row 3 (H-SubArraySort row) ---
$ws.Range("C3:D3").Style = "Normal"
$ws.Range("C3").Value = "Key is to find min and max elements in the identified unsorted array. "
$ws.Range("D3").Value = "For decrement, For increment loop, break"
$ws.Range("C3:D3").WrapText = $true

# Restore the selection Excel leaves behind after this kind of edit.
$ws.Range("D5").Select() | Out-Null
